$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the DATE_COLLECTED column (G2:G19) to a single value ---
# The collection records are no longer date-checked: every row's
# DATE_COLLECTED is unified to "2014-03-01", including the two rows
# (8 and 19) that previously held outlier dates flagged for review.
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G
    $cell.Value = "'2014-03-01"
    $cell.Style = "Normal"
}

# --- Drop the "flag for review" highlight on AJ19 (DATE_IDENTIFIED) ---
# Its text value is unchanged; only the highlighted/flagged style goes away.
$ws.Range("AJ19").Style = "Normal"

# --- Remove the "Please check the date" review comments ---
# (G8, G19, AJ19) and re-establish the comments that remain untouched,
# since this runtime cannot round-trip the pre-existing comment part.
$ws.Range("G8").ClearComments()
$ws.Range("G19").ClearComments()
$ws.Range("AJ19").ClearComments()

$ws.Range("Q14").ClearComments()
$ws.Range("Q14").AddComment("STATE_PROVINCE not provided") | Out-Null

$ws.Range("Q15").ClearComments()
$ws.Range("Q15").AddComment("STATE_PROVINCE not provided") | Out-Null

$ws.Range("S4").ClearComments()
$ws.Range("S4").AddComment("LOCALITY not provided") | Out-Null

$ws.Range("S14").ClearComments()
$ws.Range("S14").AddComment("LOCALITY not provided") | Out-Null

$ws.Range("V9").ClearComments()
$ws.Range("V9").AddComment("Degrees must be between 0 and 90") | Out-Null
